# AMOS B02 - Scrum and AMOS: update table style references.
#
# The deck's tables were pointing at the old default table style
# {C8671AF5-1898-4181-9D50-357DC450F274}; switch them to the new
# style {2E9E8560-8F3E-4D38-A19E-F5F1590397B1}.

$p = $ppt.ActivePresentation

$oldStyleId = "{C8671AF5-1898-4181-9D50-357DC450F274}"
$newStyleId = "{2E9E8560-8F3E-4D38-A19E-F5F1590397B1}"

for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($j = 1; $j -le $slide.Shapes.Count; $j++) {
        $shape = $slide.Shapes.Item($j)
        if ($shape.HasTable) {
            $table = $shape.Table
            if ($table.Style.Name -eq $oldStyleId) {
                $table.ApplyStyle($newStyleId)
            }
        }
    }
}
